$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns G..K are brand new (the sheet used to stop at F). Give them the
# same header formatting (bold font + border + centered) already used by
# the existing header cells, by copying F1's format across.
$ws.Range("F1").Copy()
$ws.Range("G1:K1").PasteSpecial(-4122)

# New header labels for columns F..K
$ws.Cells.Item(1, 6).Value = "frequency"
$ws.Cells.Item(1, 7).Value = "frequency_occurrence"
$ws.Cells.Item(1, 8).Value = "frequency_occurrence_probab"
$ws.Cells.Item(1, 9).Value = "max_probab"
$ws.Cells.Item(1, 10).Value = "max_probab_percentage"
$ws.Cells.Item(1, 11).Value = "recommended_level"

# Per-row frequency / probability data (F..K), keyed by row number.
# Rows not present here (5-8, 36, 39-82) are left blank, matching the source data.
# J is a percentage rendered as text (e.g. "75.00"); a leading apostrophe
# forces it to be stored as text instead of being coerced to the number 75.
$rowData = @{
    2 = @{ F=12; G='{"L2":9,"L1":3}'; H='{"L2":0.75,"L1":0.25}'; I=0.75; J="'75.00"; K='L2' }
    3 = @{ F=12; G='{"L2":9,"L1":3}'; H='{"L2":0.75,"L1":0.25}'; I=0.75; J="'75.00"; K='L2' }
    4 = @{ F=9; G='{"L3":9}'; H='{"L3":1.0}'; I=1; J="'100.00"; K='L3' }
    9 = @{ F=8; G='{"L2":5,"L3":3}'; H='{"L2":0.625,"L3":0.375}'; I=0.625; J="'62.50"; K='L3' }
    10 = @{ F=8; G='{"L2":8}'; H='{"L2":1.0}'; I=1; J="'100.00"; K='L2' }
    11 = @{ F=8; G='{"L2":7,"L1":1}'; H='{"L2":0.875,"L1":0.125}'; I=0.875; J="'87.50"; K='L2' }
    12 = @{ F=8; G='{"L1":4,"L2":4}'; H='{"L1":0.5,"L2":0.5}'; I=0.5; J="'50.00"; K='L2' }
    13 = @{ F=8; G='{"L3":6,"L2":2}'; H='{"L3":0.75,"L2":0.25}'; I=0.75; J="'75.00"; K='L3' }
    14 = @{ F=8; G='{"L1":8}'; H='{"L1":1.0}'; I=1; J="'100.00"; K='L1' }
    15 = @{ F=8; G='{"L3":6,"L2":2}'; H='{"L3":0.75,"L2":0.25}'; I=0.75; J="'75.00"; K='L3' }
    16 = @{ F=8; G='{"L3":7,"L2":1}'; H='{"L3":0.875,"L2":0.125}'; I=0.875; J="'87.50"; K='L3' }
    17 = @{ F=8; G='{"L3":7,"L2":1}'; H='{"L3":0.875,"L2":0.125}'; I=0.875; J="'87.50"; K='L3' }
    18 = @{ F=8; G='{"L3":8}'; H='{"L3":1.0}'; I=1; J="'100.00"; K='L3' }
    19 = @{ F=8; G='{"L3":4,"L2":4}'; H='{"L3":0.5,"L2":0.5}'; I=0.5; J="'50.00"; K='L3' }
    20 = @{ F=8; G='{"L3":6,"L2":2}'; H='{"L3":0.75,"L2":0.25}'; I=0.75; J="'75.00"; K='L3' }
    21 = @{ F=8; G='{"L3":5,"L2":3}'; H='{"L3":0.625,"L2":0.375}'; I=0.625; J="'62.50"; K='L3' }
    22 = @{ F=8; G='{"L3":6,"L1":1,"L2":1}'; H='{"L3":0.75,"L1":0.125,"L2":0.125}'; I=0.75; J="'75.00"; K='L3' }
    23 = @{ F=8; G='{"L3":8}'; H='{"L3":1.0}'; I=1; J="'100.00"; K='L3' }
    24 = @{ F=8; G='{"L3":8}'; H='{"L3":1.0}'; I=1; J="'100.00"; K='L3' }
    25 = @{ F=8; G='{"L3":8}'; H='{"L3":1.0}'; I=1; J="'100.00"; K='L3' }
    26 = @{ F=8; G='{"L3":8}'; H='{"L3":1.0}'; I=1; J="'100.00"; K='L3' }
    27 = @{ F=14; G='{"L3":10,"L2":3,"L1":1}'; H='{"L3":0.7142857143,"L2":0.2142857143,"L1":0.0714285714}'; I=0.7142857142857143; J="'71.43"; K='L3' }
    28 = @{ F=14; G='{"L3":10,"L2":3,"L1":1}'; H='{"L3":0.7142857143,"L2":0.2142857143,"L1":0.0714285714}'; I=0.7142857142857143; J="'71.43"; K='L3' }
    29 = @{ F=8; G='{"L1":8}'; H='{"L1":1.0}'; I=1; J="'100.00"; K='L1' }
    30 = @{ F=8; G='{"L3":8}'; H='{"L3":1.0}'; I=1; J="'100.00"; K='L3' }
    31 = @{ F=8; G='{"L3":8}'; H='{"L3":1.0}'; I=1; J="'100.00"; K='L3' }
    32 = @{ F=8; G='{"L3":7,"L2":1}'; H='{"L3":0.875,"L2":0.125}'; I=0.875; J="'87.50"; K='L3' }
    33 = @{ F=8; G='{"L3":5,"L2":2,"L1":1}'; H='{"L3":0.625,"L2":0.25,"L1":0.125}'; I=0.625; J="'62.50"; K='L3' }
    34 = @{ F=8; G='{"L2":5,"L1":2,"L3":1}'; H='{"L2":0.625,"L1":0.25,"L3":0.125}'; I=0.625; J="'62.50"; K='L3' }
    35 = @{ F=9; G='{"L3":5,"L2":4}'; H='{"L3":0.5555555556,"L2":0.4444444444}'; I=0.5555555555555556; J="'55.56"; K='L3' }
    37 = @{ F=8; G='{"L3":5,"L2":3}'; H='{"L3":0.625,"L2":0.375}'; I=0.625; J="'62.50"; K='L3' }
    38 = @{ F=8; G='{"L1":5,"L2":3}'; H='{"L1":0.625,"L2":0.375}'; I=0.625; J="'62.50"; K='L2' }
}

for ($r = 2; $r -le 82; $r++) {
    if ($rowData.ContainsKey($r)) {
        $d = $rowData[$r]
        $ws.Cells.Item($r, 6).Value = $d.F
        $ws.Cells.Item($r, 7).Value = $d.G
        $ws.Cells.Item($r, 8).Value = $d.H
        $ws.Cells.Item($r, 9).Value = $d.I
        $ws.Cells.Item($r, 10).Value = $d.J
        $ws.Cells.Item($r, 11).Value = $d.K
    } else {
        $ws.Cells.Item($r, 6).Value = ""
        $ws.Cells.Item($r, 7).Value = ""
        $ws.Cells.Item($r, 8).Value = ""
        $ws.Cells.Item($r, 9).Value = ""
        $ws.Cells.Item($r, 10).Value = ""
        $ws.Cells.Item($r, 11).Value = ""
    }
}
